$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("run_1")
$ws.Cells.Item(2, 6).Value = 30.23210453987122
$ws.Cells.Item(3, 6).Value = 29.65205216407776
$ws.Cells.Item(4, 6).Value = 29.72540855407715
$ws.Cells.Item(5, 6).Value = 29.6457314491272
$ws.Cells.Item(6, 6).Value = 29.95413899421692
$ws.Cells.Item(7, 6).Value = 30.03920912742615
$ws.Cells.Item(8, 6).Value = 29.79284834861756
$ws.Cells.Item(9, 6).Value = 29.85354518890381
$ws.Cells.Item(10, 6).Value = 29.71895694732666
$ws.Cells.Item(11, 6).Value = 30.161936044693
$ws.Cells.Item(12, 6).Value = 29.85380482673645
$ws.Cells.Item(13, 6).Value = 29.6971070766449
$ws.Cells.Item(14, 6).Value = 29.79974770545959
$ws.Cells.Item(15, 6).Value = 29.88181185722351
$ws.Cells.Item(16, 6).Value = 29.89919328689575
$ws.Cells.Item(17, 6).Value = 29.8801155090332
$ws.Cells.Item(18, 6).Value = 29.81754326820373
$ws.Cells.Item(19, 6).Value = 29.78310918807984
$ws.Cells.Item(20, 6).Value = 29.93346691131592
$ws.Cells.Item(21, 6).Value = 30.20338106155396

$ws = $wb.Worksheets.Item("run_2")
$ws.Cells.Item(2, 6).Value = 30.2936327457428
$ws.Cells.Item(3, 6).Value = 29.76337814331055
$ws.Cells.Item(4, 6).Value = 29.82246732711792
$ws.Cells.Item(5, 6).Value = 29.64490723609924
$ws.Cells.Item(6, 6).Value = 29.84499406814575
$ws.Cells.Item(7, 6).Value = 29.66664791107178
$ws.Cells.Item(8, 6).Value = 29.86287498474121
$ws.Cells.Item(9, 6).Value = 29.66117691993713
$ws.Cells.Item(10, 6).Value = 29.88728880882263
$ws.Cells.Item(11, 6).Value = 29.92850255966187
$ws.Cells.Item(12, 6).Value = 29.71263527870178
$ws.Cells.Item(13, 6).Value = 29.60158276557922
$ws.Cells.Item(14, 6).Value = 29.68245768547058
$ws.Cells.Item(15, 6).Value = 29.65557599067688
$ws.Cells.Item(16, 6).Value = 29.66921329498291
$ws.Cells.Item(17, 6).Value = 29.84323382377625
$ws.Cells.Item(18, 6).Value = 29.9485890865326
$ws.Cells.Item(19, 6).Value = 29.89292907714844
$ws.Cells.Item(20, 6).Value = 29.92488503456116
$ws.Cells.Item(21, 6).Value = 30.16585803031921

$ws = $wb.Worksheets.Item("run_3")
$ws.Cells.Item(2, 6).Value = 30.25470042228699
$ws.Cells.Item(3, 6).Value = 29.97967147827148
$ws.Cells.Item(4, 6).Value = 29.66020727157593
$ws.Cells.Item(5, 6).Value = 29.62090754508972
$ws.Cells.Item(6, 6).Value = 29.68105697631836
$ws.Cells.Item(7, 6).Value = 29.84277629852295
$ws.Cells.Item(8, 6).Value = 29.92955422401428
$ws.Cells.Item(9, 6).Value = 29.77686142921448
$ws.Cells.Item(10, 6).Value = 29.67220044136048
$ws.Cells.Item(11, 6).Value = 29.9480185508728
$ws.Cells.Item(12, 6).Value = 29.84974384307861
$ws.Cells.Item(13, 6).Value = 29.65391945838928
$ws.Cells.Item(14, 6).Value = 29.75651860237122
$ws.Cells.Item(15, 6).Value = 29.65123081207276
$ws.Cells.Item(16, 6).Value = 29.756352186203
$ws.Cells.Item(17, 6).Value = 29.5340564250946
$ws.Cells.Item(18, 6).Value = 29.72342681884766
$ws.Cells.Item(19, 6).Value = 29.78702688217163
$ws.Cells.Item(20, 6).Value = 29.68014740943909
$ws.Cells.Item(21, 6).Value = 29.81495952606201

$ws = $wb.Worksheets.Item("run_4")
$ws.Cells.Item(2, 6).Value = 30.07565569877625
$ws.Cells.Item(3, 6).Value = 29.66848754882812
$ws.Cells.Item(4, 6).Value = 29.83465147018433
$ws.Cells.Item(5, 6).Value = 29.66834211349488
$ws.Cells.Item(6, 6).Value = 30.00508165359497
$ws.Cells.Item(7, 6).Value = 29.82563090324402
$ws.Cells.Item(8, 6).Value = 29.87941312789917
$ws.Cells.Item(9, 6).Value = 29.7331235408783
$ws.Cells.Item(10, 6).Value = 29.84268403053284
$ws.Cells.Item(11, 6).Value = 30.00035500526428
$ws.Cells.Item(12, 6).Value = 29.91219663619995
$ws.Cells.Item(13, 6).Value = 29.87175893783569
$ws.Cells.Item(14, 6).Value = 29.96086812019348
$ws.Cells.Item(15, 6).Value = 29.88192677497864
$ws.Cells.Item(16, 6).Value = 29.85164546966553
$ws.Cells.Item(17, 6).Value = 29.50654721260071
$ws.Cells.Item(18, 6).Value = 29.85675668716431
$ws.Cells.Item(19, 6).Value = 29.82875990867615
$ws.Cells.Item(20, 6).Value = 29.62727069854736
$ws.Cells.Item(21, 6).Value = 30.08136391639709

$ws = $wb.Worksheets.Item("run_5")
$ws.Cells.Item(2, 6).Value = 30.30006098747253
$ws.Cells.Item(3, 6).Value = 30.76500272750854
$ws.Cells.Item(4, 6).Value = 30.83791518211365
$ws.Cells.Item(5, 6).Value = 30.52731418609619
$ws.Cells.Item(6, 6).Value = 30.65078496932984
$ws.Cells.Item(7, 6).Value = 30.49234676361084
$ws.Cells.Item(8, 6).Value = 30.81153440475464
$ws.Cells.Item(9, 6).Value = 30.13240361213684
$ws.Cells.Item(10, 6).Value = 29.94280290603638
$ws.Cells.Item(11, 6).Value = 30.03670763969421
$ws.Cells.Item(12, 6).Value = 29.84545993804932
$ws.Cells.Item(13, 6).Value = 29.74356484413147
$ws.Cells.Item(14, 6).Value = 29.88382577896118
$ws.Cells.Item(15, 6).Value = 29.63089942932129
$ws.Cells.Item(16, 6).Value = 29.82643055915833
$ws.Cells.Item(17, 6).Value = 29.81170177459717
$ws.Cells.Item(18, 6).Value = 29.80436182022095
$ws.Cells.Item(19, 6).Value = 29.78512978553772
$ws.Cells.Item(20, 6).Value = 29.87981796264648
$ws.Cells.Item(21, 6).Value = 30.15430998802185

Write-Output "Updated Epoch Time (column F) for run_1..run_5, rows 2-21"